# ---------------------------------------------------------------------------
# 9-month Size workbook update
#   * Adds a note to K3 explaining that some records were added with
#     size = 2mm
#   * Appends 8 new data rows (324-331) for the HL6-AMB / HL6-LOW cohorts
#     that were added to the dataset as size = 2mm
#   * Nudges column K a little wider so the new note is readable
#   * Restores the selection to D5 (no more frozen/scrolled topLeftCell)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New note in K3 --------------------------------------------------------
$ws.Range("K3").Value = "These were added to the dataset as size = 2mm "

# --- Give column K a bit more breathing room for the note -----------------
$ws.Columns.Item(11).ColumnWidth = 13.25

# --- New row 324: HL6-AMB, all measurements = 2mm --------------------------
$ws.Range("A324").Value = "NA"
$ws.Range("B324").Value = "HL6-AMB"
$ws.Range("C324").Value = "HL"
$ws.Range("D324").Value = "6-AMB"
$ws.Range("E324").Value = "NA"
$ws.Range("F324").Value = "NA"
$ws.Range("G324:T324").Value = 2

# --- New rows 325-331: HL6-LOW, all measurements = 2mm ---------------------
for ($r = 325; $r -le 331; $r++) {
    $ws.Range("A$r").Value = "NA"
    $ws.Range("B$r").Value = "HL6-LOW"
    $ws.Range("C$r").Value = "HL"
    $ws.Range("D$r").Value = "6-LOW"
    $ws.Range("E$r").Value = "NA"
    $ws.Range("F$r").Value = "NA"
    $ws.Range("G$r" + ":AE$r").Value = 2
}

# --- Move the saved selection/cursor to D5 (also clears topLeftCell) ------
$ws.Range("D5").Select()
